$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to stay as text
# (matching the original inline-string/text formatting), so we set NumberFormat to
# "@" (Text) before assigning them. Cells whose values are not number-like (they
# contain letters, URLs, percent signs, padding spaces or multiple dots) are safe
# to assign directly.

$ws.Range("D2").Value = '22.386.09'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.570.74'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.67'
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3760'
$ws.Range("E7").Value = '  +2.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.70'
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3409'
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07616'
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.13'
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.007'
$ws.Range("E14").Value = '  -0.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.958'
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").Value = '1.571.51'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.16'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06757'
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.72'
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.192'
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").Value = '22.389.09'
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.389'
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.679'
$ws.Range("E26").Value = '  -8.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.09'
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.44'
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.032'
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.56'
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("D31").Value = '1.746.36'
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.016'
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.002'
$ws.Range("E33").Value = '  -4.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.089'
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '10.13'
$ws.Range("E35").Value = '  -0.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08479'
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02539'
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.385'
$ws.Range("E38").Value = '  +10.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2305'
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06505'
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.424'
$ws.Range("E41").Value = '  -2.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.39'
$ws.Range("E42").Value = '  -2.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6330'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.98'
$ws.Range("E45").Value = '  -3.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.806'
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5943'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.084'
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.280'
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.42'
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07317'
$ws.Range("E51").Value = '  +0.29%  '
